# Generate Report for Archive
#
# 1. Replace every occurrence of the shared string "Ready for handoff"
#    with "In Translation" (Overview!E2:F4 and Status column C2:C4 on
#    both the "zh-cn" and "de-de" sheets).
# 2. Narrow the "zh-cn"/"de-de" status columns: Overview columns E & F,
#    and column C on the "zh-cn"/"de-de" sheets, from ~17.216 to
#    ~13.410 characters wide.

$wb = $excel.ActiveWorkbook

# ---- 1. "Ready for handoff" -> "In Translation" ------------------------
# Overview: status columns E (zh-cn) and F (de-de), rows 2-4.
$overview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 4; $r++) {
    foreach ($c in 5, 6) {
        $cell = $overview.Cells.Item($r, $c)
        if ($cell.Text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# zh-cn / de-de: "Status" column C, rows 2-4.
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
foreach ($ws in $zhcn, $dede) {
    for ($r = 2; $r -le 4; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        if ($cell.Text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# ---- 2. Narrow the status columns from ~17.22 to ~13.41 wide -----------
$overview.Range("E:E").ColumnWidth = 13.4101845877511
$overview.Range("F:F").ColumnWidth = 13.4101845877511
$zhcn.Range("C:C").ColumnWidth = 13.4101845877511
$dede.Range("C:C").ColumnWidth = 13.4101845877511
